$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the _GoBack bookmark that sits after "2023"
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $goBack = $d.Bookmarks.Item("_GoBack")
    $goBack.Delete()
}

# ------------------------------------------------------------------
# 2) Before the "NBA tem uma longa história..." paragraph, add a new
#    paragraph "-- falar sobre o basquete e a nba --" and prefix the
#    existing paragraph text with "A ".
# ------------------------------------------------------------------
$contextParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith("NBA tem uma longa hist")) {
        $contextParaIndex = $i
        break
    }
}

$contextPara = $d.Paragraphs($contextParaIndex)
$null = $contextPara.Range.InsertParagraphBefore()

$introPara = $d.Paragraphs($contextParaIndex)
$introRange = $introPara.Range
$introRange.Collapse(1)
$introXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="708"/></w:pPr><w:r><w:t xml:space="preserve">-- falar sobre o basquete e a </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>nba</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> --</w:t></w:r></w:p>'
$null = $introRange.InsertXML($introXml)

$contextParaIndex = $contextParaIndex + 1
$contextPara = $d.Paragraphs($contextParaIndex)
$contextRange = $contextPara.Range
$contextRange.Collapse(1)
$null = $contextRange.InsertBefore("A ")

# ------------------------------------------------------------------
# 3) Split the "Justificativa" paragraph ("<tab>Eu sou bastante fã...")
#    into: "<tab>--  explicar melhor o motivo de eu gostar do esporte
#    e como que eu conheci --" (with the _GoBack bookmark at the end)
#    followed by its own paragraph "Eu sou bastante fã de basquete..."
# ------------------------------------------------------------------
$justParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t.IndexOf("Eu sou bastante f") -ge 0) {
        $justParaIndex = $i
        break
    }
}

$justPara = $d.Paragraphs($justParaIndex)
$justRange = $justPara.Range
$justRange.Collapse(1)
$justXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00A72CB7" w:rsidRDefault="00A72CB7" w:rsidP="00A72CB7"><w:r><w:tab/></w:r><w:r><w:t>--</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">  </w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>explicar melhor o motivo de eu gostar do esporte e como que eu conheci --</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="708"/></w:pPr><w:r><w:t>Eu sou bastante fã de basquete e NBA, portanto fico bastante interessado em prover para pessoas informações sobre o esporte que amo.</w:t></w:r></w:p>'
$null = $justRange.InsertXML($justXml)
